$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits after the
#    "Rigidbody ,gravity -> iskinematic" run.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Turn two of the blank paragraphs that follow
#    "smoke effect when low health?" into the new "Animation:" / setKey
#    notes, and re-create the "_GoBack" bookmark at the end of the new
#    content (this is where Word's cursor/_GoBack position moved to after
#    the edit).
# ---------------------------------------------------------------------------

# Find the "smoke effect..." paragraph as an anchor.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "smoke effect when low health") {
        $anchorIndex = $i
        break
    }
}

# Walk forward through the run of empty paragraphs that follows it and
# grab the 4th and 5th ones (1-indexed) - these are the two paragraphs
# that get replaced with real content.
$emptyIdx = @()
$j = $anchorIndex + 1
while ($emptyIdx.Count -lt 9 -and $d.Paragraphs.Item($j).Range.Text.Trim().Length -eq 0) {
    $emptyIdx += $j
    $j = $j + 1
}

$targetStartIndex = $emptyIdx[3]
$targetEndIndex = $emptyIdx[4]

$p1 = $d.Paragraphs.Item($targetStartIndex)
$p2 = $d.Paragraphs.Item($targetEndIndex)
$r = $d.Range($p1.Range.Start, $p2.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr>
<w:r><w:rPr><w:color w:val="FF0000"/><w:highlight w:val="yellow"/></w:rPr><w:t>Animation:</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>设置</w:t></w:r>
<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>,setKey!</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r.InsertXML($xml) | Out-Null
